# Update min_price (D) and autonova-nkz.ru_price (F) for the rows whose
# scraped price changed. Column F stores the price as text (inlineStr in
# the source workbook), so it is explicitly formatted as Text before the
# value is written to avoid Excel auto-converting the numeric-looking
# string back into a Number cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D107").Value = 291800
$ws.Range("F107").NumberFormat = "@"
$ws.Range("F107").Value = "291800"

$ws.Range("D108").Value = 286300
$ws.Range("F108").NumberFormat = "@"
$ws.Range("F108").Value = "286300"

$ws.Range("D109").Value = 293300
$ws.Range("F109").NumberFormat = "@"
$ws.Range("F109").Value = "293300"

$ws.Range("D110").Value = 280400
$ws.Range("F110").NumberFormat = "@"
$ws.Range("F110").Value = "280400"

$ws.Range("D111").Value = 267900
$ws.Range("F111").NumberFormat = "@"
$ws.Range("F111").Value = "267900"

$ws.Range("D112").Value = 300300
$ws.Range("F112").NumberFormat = "@"
$ws.Range("F112").Value = "300300"

$ws.Range("D120").Value = 460500
$ws.Range("F120").NumberFormat = "@"
$ws.Range("F120").Value = "460500"

$ws.Range("D122").Value = 424000
$ws.Range("F122").NumberFormat = "@"
$ws.Range("F122").Value = "424000"

$ws.Range("D123").Value = 449000
$ws.Range("F123").NumberFormat = "@"
$ws.Range("F123").Value = "449000"

$ws.Range("D126").Value = 396900
$ws.Range("F126").NumberFormat = "@"
$ws.Range("F126").Value = "396900"

$ws.Range("D127").Value = 481900
$ws.Range("F127").NumberFormat = "@"
$ws.Range("F127").Value = "481900"

$ws.Range("D129").Value = 565900
$ws.Range("F129").NumberFormat = "@"
$ws.Range("F129").Value = "565900"

$ws.Range("D130").Value = 404300
$ws.Range("F130").NumberFormat = "@"
$ws.Range("F130").Value = "404300"

$ws.Range("D131").Value = 395300
$ws.Range("F131").NumberFormat = "@"
$ws.Range("F131").Value = "395300"

